$wb = $excel.ActiveWorkbook
$ws3 = $wb.Worksheets.Item("Alunni")
$ws3.Rows.Item(5).Delete()
$ws3.Rows.Item(8).Delete()
$ws3.Rows.Item(11).Delete()
$ws3.Rows.Item(14).Delete()
$ws3.Rows.Item(17).Delete()
$ws3.Rows.Item(20).Delete()
$ws3.Rows.Item(23).Delete()
$ws3.Rows.Item(26).Delete()
Write-Host "done"
